# Fix missing length (column D) and weight (column F) values in the
# rebar schedule table on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> @{ D = length; F = weight (optional) }
$updates = @{
    3  = @{ D = 3000 }
    4  = @{ D = 3000; F = 19212 }
    5  = @{ D = 571;  F = 3657 }
    6  = @{ D = 1400; F = 8966 }
    7  = @{ D = 2200; F = 14089 }
    8  = @{ D = 510;  F = 3269 }
    9  = @{ D = 571;  F = 3657 }
    10 = @{ D = 2000; F = 12808 }
    11 = @{ D = 300;  F = 2372 }
    12 = @{ D = 300;  F = 1508 }
    13 = @{ D = 300;  F = 1192 }
    14 = @{ D = 300;  F = 913 }
    15 = @{ D = 300;  F = 670 }
    16 = @{ D = 300;  F = 466 }
    17 = @{ D = 300;  F = 299 }
    18 = @{ D = 300;  F = 1921 }
    19 = @{ D = 300;  F = 168 }
    20 = @{ D = 700;  F = 5534 }
    21 = @{ D = 500;  F = 3953 }
    22 = @{ D = 1060; F = 8380 }
    23 = @{ D = 560;  F = 4427 }
    24 = @{ D = 600;  F = 4744 }
    25 = @{ D = 1120; F = 8855 }
    26 = @{ D = 620;  F = 4902 }
    27 = @{ D = 700;  F = 5534 }
    28 = @{ D = 700;  F = 5534 }
    29 = @{ D = 700;  F = 5534 }
    30 = @{ D = 760;  F = 6009 }
    31 = @{ D = 660;  F = 5218 }
    32 = @{ D = 560;  F = 4427 }
    33 = @{ D = 745;  F = 3744 }
    34 = @{ D = 745;  F = 3744 }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey('D')) {
        $ws.Cells.Item($row, 4).Value = $vals['D']
    }
    if ($vals.ContainsKey('F')) {
        $ws.Cells.Item($row, 6).Value = $vals['F']
    }
}
